$d = $word.ActiveDocument

function Get-ParaIndexByStart($doc, $startPos) {
    $idx = 1
    foreach ($para in $doc.Paragraphs) {
        if ($para.Range.Start -eq $startPos) {
            return $idx
        }
        $idx++
    }
    return -1
}

# --- Insert two new bullet paragraphs right after the "Streamlined
# communication..." bullet (before "Consulted with multiple groups...") in
# the STRATCOM Strategic Plans Officer / Security Manager job block.
$anchor1 = $d.Range(0, 0)
$anchor1.Find.Execute("Streamlined communication up chain of command", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$anchorIdx1 = Get-ParaIndexByStart $d $anchor1.Start
$anchorPara1 = $d.Paragraphs($anchorIdx1)
$anchorPara1.Range.InsertParagraphAfter()
$newPara1 = $d.Paragraphs($anchorIdx1 + 1)
$newPara1.Range.InsertXML('<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Compact"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1013"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Worked with and around multiple applications on a regular basis such as Global Strike Planning Aid (GSPA, now named ACME), Joint Targetting Toolbox (JTT), Digital Imagery Exploitation Engine (DIEE), SPA, and IMEA</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Compact"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1013"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Qualified as a Collateral Damage Estimate (CDE) analyst</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

# --- Insert a new certification bullet right after the "CERTIFICATIONS"
# heading (before the "2015 ... CISSP" bullet).
$anchor2 = $d.Range(0, 0)
$anchor2.Find.Execute("CERTIFICATIONS", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$anchorIdx2 = Get-ParaIndexByStart $d $anchor2.Start
$anchorPara2 = $d.Paragraphs($anchorIdx2)
$anchorPara2.Range.InsertParagraphAfter()
$newPara2 = $d.Paragraphs($anchorIdx2 + 1)
$newPara2.Range.InsertXML('<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Compact"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1017"/></w:numPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">????</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">Offensive Security Certified Professional (OSCP)</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

Write-Host "Edits applied."
